# Apply the new table style ({D32551E2-C605-4F60-801C-5EE43A18CBF1}) to the
# three tables that currently use the custom "Table_0" style
# ({916EA1F1-D4A7-4BC2-A8A2-E3E5C7E99520}). Each of these tables is the
# first (and, on slide 16, only) shape on its slide.
$p = $ppt.ActivePresentation

$newStyleId = "{D32551E2-C605-4F60-801C-5EE43A18CBF1}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId, $true)
    }
}
